$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook calculation: switch off manual calc mode (calcMode="manual" removed) ---
$excel.Calculation = -4105   # xlCalculationAutomatic

# --- Row 17: repurpose existing last row from "Any/waitForPageToRender/..." ---
# to "Text Field/textBoxShouldHaveValue/..."
$ws.Range("B17").Value = "Text Field"
$ws.Range("C17").Value = "textBoxShouldHaveValue"
$ws.Range("C17").WrapText = $true
$ws.Range("D17").Value = "Accepts two parameters @locator and @testData. It gets the text from textBox and validates against the @testData provided. If the validation fails testing should still continue"
$ws.Range("D17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 29

# --- Row 18: re-create the row that used to be row 17 (waitForPageToRender) ---
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Any"
$ws.Range("C18").Value = "waitForPageToRender"
$ws.Range("D18").Value = "Accepts no parameters and waits for page to render itself. It does that by observing network state. Step never triggers test failure and waits for maximum period of 30 seconds"
$ws.Range("D18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 29

# --- Row 19: elementShouldNotBePresent ---
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Any"
$ws.Range("C19").Value = "elementShouldNotBePresent"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("D19").Value = "Accepts no parameters and verifies element is not available in DOM. Returns true if element is not available in DOM"
$ws.Range("D19").WrapText = $true
$ws.Rows.Item(19).RowHeight = 29

# --- Rows 20-21: Button / isButtonEnabled / isButtonDisabled ---
# Short-label columns (B,C) are filled first for both rows, then the longer
# description column (D) is filled for both rows, matching the order in
# which the underlying shared-string table was populated.
$ws.Range("B20").Value = "Button"
$ws.Range("C20").Value = "isButtonEnabled"
$ws.Range("B21").Value = "Button"
$ws.Range("C21").Value = "isButtonDisabled"

$ws.Range("D20").Value = "Accepts one parameter @locator. Checks if the Button is enabled. Returns true if the button is enabled"
$ws.Range("D21").Value = "Accepts one parameter @locator. Checks if the Button is disabled Returns true if the button is disabled"

$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20

$ws.Range("D20").WrapText = $true
$ws.Range("D21").WrapText = $true
$ws.Rows.Item(20).RowHeight = 29
$ws.Rows.Item(21).RowHeight = 29

# --- Selection / scroll position ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D21").Select()
